# Update the "G column" formulas on the "老師上課 + 監考時數" sheet.
# Rows 2..167 (the specific data rows listed below) previously repeated a
# long SUM-based expression; that full computation now lives only in G168,
# and each data row simply multiplies it by its own C-column ratio.
# G168 itself is updated to round the result to a whole number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("老師上課 + 監考時數")

$rows = @(2,7,11,14,16,18,20,24,26,27,29,33,36,37,40,43,44,46,50,54,57,59,64,68,70,76,77,79,82,85,90,94,96,99,102,105,108,110,112,114,116,118,120,121,122,123,125,130,132,133,135,138,142,145,149,151,154,158,161,162,163,164,166,167)

foreach ($r in $rows) {
    $ws.Range("G$r").Formula = "=G168*C$r"
}

$ws.Range("G168").Formula = "=ROUND((SUM(`$D`$2:`$D`$167*`$C`$2:`$C`$167)-F161-F162-F163-F164+SUM(`$E`$2:`$E`$167))/SUM(`$C`$2:`$C`$167),0)"
